# Update "想去人数" (number of interested attendees) counts in column F
# across the four worksheets, reflecting the latest scrape at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1886
$ws.Range("F3").Value = 1516
$ws.Range("F4").Value = 883
$ws.Range("F5").Value = 768
$ws.Range("F6").Value = 13333
$ws.Range("F7").Value = 13197
$ws.Range("F11").Value = 558
$ws.Range("F13").Value = 677
$ws.Range("F17").Value = 74
$ws.Range("F20").Value = 257
$ws.Range("F21").Value = 289
$ws.Range("F22").Value = 423
$ws.Range("F23").Value = 759
$ws.Range("F24").Value = 17

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 94
$ws.Range("F7").Value = 123

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 42

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1886
$ws.Range("F4").Value = 1516
$ws.Range("F5").Value = 883
$ws.Range("F6").Value = 94
$ws.Range("F7").Value = 768
$ws.Range("F8").Value = 13333
$ws.Range("F9").Value = 13197
$ws.Range("F13").Value = 558
$ws.Range("F15").Value = 677
$ws.Range("F21").Value = 74
$ws.Range("F25").Value = 42
$ws.Range("F27").Value = 257
$ws.Range("F28").Value = 289
$ws.Range("F29").Value = 423
$ws.Range("F30").Value = 759
$ws.Range("F31").Value = 123
$ws.Range("F33").Value = 17

$wb.Save()
